$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row: marks awarded per correct answer
$ws.Range("B11").Value = 5

# Update "Total" row: total score achieved and score/max text
$ws.Range("B12").Value = 130
$ws.Range("E12").Value = "130/140"
